$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151; existing rows 151-171 shift down to 152-172.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new weekly record.
$ws.Cells.Item(151, 1).Value = 9
$ws.Cells.Item(151, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(151, 3).Value = "Metropolitana"
$ws.Cells.Item(151, 4).Value = 44984
$ws.Cells.Item(151, 5).Value = 13
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100103
$ws.Cells.Item(151, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(151, 9).Value = 100103002
$ws.Cells.Item(151, 10).Value = "Ciruela"
$ws.Cells.Item(151, 11).Value = "Friar"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 12
$ws.Cells.Item(151, 14).Value = 180000
$ws.Cells.Item(151, 15).Value = 180000
$ws.Cells.Item(151, 16).Value = 180000
$ws.Cells.Item(151, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(151, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(151, 19).Value = 400
$ws.Cells.Item(151, 20).Value = 450
